# Update the "Förändrad" (changed) date column (C) for every existing data
# row (2..392) from 2023-09-20 (45189) to 2023-09-21 (45190) in one shot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C392").Value2 = 45190

# Row 392 picks up an explicit (default) row height, matching the diff.
$ws.Rows.Item(392).RowHeight = 15

# New row 393: "A 44495-2023"
$ws.Range("A393").Value2 = "A 44495-2023"
$ws.Range("B393").Value2 = 45189
$ws.Range("B393").NumberFormat = "YYYY-MM-DD"
$ws.Range("C393").Value2 = 45190
$ws.Range("C393").NumberFormat = "YYYY-MM-DD"
$ws.Range("D393").Value2 = "VÄSTERBOTTENS LÄN"
$ws.Range("E393").Value2 = "MALÅ"
$ws.Range("G393").Value2 = 1.1
$ws.Range("H393:Q393").Value2 = 0
$ws.Range("R393").Value2 = ""
$ws.Range("R393").WrapText = $true
$ws.Rows.Item(393).RowHeight = 15

# New row 394: "A 44596-2023"
$ws.Range("A394").Value2 = "A 44596-2023"
$ws.Range("B394").Value2 = 45189
$ws.Range("B394").NumberFormat = "YYYY-MM-DD"
$ws.Range("C394").Value2 = 45190
$ws.Range("C394").NumberFormat = "YYYY-MM-DD"
$ws.Range("D394").Value2 = "VÄSTERBOTTENS LÄN"
$ws.Range("E394").Value2 = "MALÅ"
$ws.Range("G394").Value2 = 4.1
$ws.Range("H394:Q394").Value2 = 0
$ws.Range("R394").Value2 = ""
$ws.Range("R394").WrapText = $true
